$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
${ws}.Cells.Item(2, 4).NumberFormat = "@"
${ws}.Cells.Item(2, 4).Value = "27.532.50"
${ws}.Cells.Item(2, 5).NumberFormat = "@"
${ws}.Cells.Item(2, 5).Value = "  +5.29%  "

# Row 3
${ws}.Cells.Item(3, 4).NumberFormat = "@"
${ws}.Cells.Item(3, 4).Value = "1.724.47"
${ws}.Cells.Item(3, 5).NumberFormat = "@"
${ws}.Cells.Item(3, 5).Value = "  +4.17%  "

# Row 4
${ws}.Cells.Item(4, 4).NumberFormat = "@"
${ws}.Cells.Item(4, 4).Value = "1.003"
${ws}.Cells.Item(4, 5).NumberFormat = "@"
${ws}.Cells.Item(4, 5).Value = "  +0.08%  "

# Row 5
${ws}.Cells.Item(5, 4).NumberFormat = "@"
${ws}.Cells.Item(5, 4).Value = "225.80"
${ws}.Cells.Item(5, 5).NumberFormat = "@"
${ws}.Cells.Item(5, 5).Value = "  +3.27%  "

# Row 6
${ws}.Cells.Item(6, 4).NumberFormat = "@"
${ws}.Cells.Item(6, 4).Value = "0.5374"
${ws}.Cells.Item(6, 5).NumberFormat = "@"
${ws}.Cells.Item(6, 5).Value = "  +2.62%  "

# Row 7
${ws}.Cells.Item(7, 5).NumberFormat = "@"
${ws}.Cells.Item(7, 5).Value = "  +0.03%  "

# Row 8
${ws}.Cells.Item(8, 4).NumberFormat = "@"
${ws}.Cells.Item(8, 4).Value = "0.2675"
${ws}.Cells.Item(8, 5).NumberFormat = "@"
${ws}.Cells.Item(8, 5).Value = "  +0.64%  "

# Row 9
${ws}.Cells.Item(9, 4).NumberFormat = "@"
${ws}.Cells.Item(9, 4).Value = "0.06604"
${ws}.Cells.Item(9, 5).NumberFormat = "@"
${ws}.Cells.Item(9, 5).Value = "  +3.95%  "

# Row 10
${ws}.Cells.Item(10, 4).NumberFormat = "@"
${ws}.Cells.Item(10, 4).Value = "21.78"
${ws}.Cells.Item(10, 5).NumberFormat = "@"
${ws}.Cells.Item(10, 5).Value = "  +5.94%  "

# Row 11
${ws}.Cells.Item(11, 4).NumberFormat = "@"
${ws}.Cells.Item(11, 4).Value = "0.07733"
${ws}.Cells.Item(11, 5).NumberFormat = "@"
${ws}.Cells.Item(11, 5).Value = "  +0.51%  "

# Row 12
${ws}.Cells.Item(12, 4).NumberFormat = "@"
${ws}.Cells.Item(12, 4).Value = "4.617"
${ws}.Cells.Item(12, 5).NumberFormat = "@"
${ws}.Cells.Item(12, 5).Value = "  +0.18%  "

# Row 13
${ws}.Cells.Item(13, 4).NumberFormat = "@"
${ws}.Cells.Item(13, 4).Value = "1.723.81"
${ws}.Cells.Item(13, 5).NumberFormat = "@"
${ws}.Cells.Item(13, 5).Value = "  +1.37%  "

# Row 14
${ws}.Cells.Item(14, 4).NumberFormat = "@"
${ws}.Cells.Item(14, 4).Value = "1.962.16"
${ws}.Cells.Item(14, 5).NumberFormat = "@"
${ws}.Cells.Item(14, 5).Value = "  +4.17%  "

# Row 15
${ws}.Cells.Item(15, 4).NumberFormat = "@"
${ws}.Cells.Item(15, 4).Value = "0.5876"
${ws}.Cells.Item(15, 5).NumberFormat = "@"
${ws}.Cells.Item(15, 5).Value = "  +4.61%  "

# Row 16
${ws}.Cells.Item(16, 4).NumberFormat = "@"
${ws}.Cells.Item(16, 4).Value = "0.0₅8303"
${ws}.Cells.Item(16, 5).NumberFormat = "@"
${ws}.Cells.Item(16, 5).Value = "  +1.35%  "

# Row 17
${ws}.Cells.Item(17, 4).NumberFormat = "@"
${ws}.Cells.Item(17, 4).Value = "68.01"
${ws}.Cells.Item(17, 5).NumberFormat = "@"
${ws}.Cells.Item(17, 5).Value = "  +3.92%  "

# Row 18
${ws}.Cells.Item(18, 4).NumberFormat = "@"
${ws}.Cells.Item(18, 4).Value = "27.553.68"
${ws}.Cells.Item(18, 5).NumberFormat = "@"
${ws}.Cells.Item(18, 5).Value = "  +5.45%  "

# Row 19
${ws}.Cells.Item(19, 4).NumberFormat = "@"
${ws}.Cells.Item(19, 4).Value = "221.77"
${ws}.Cells.Item(19, 5).NumberFormat = "@"
${ws}.Cells.Item(19, 5).Value = "  +15.31%  "

# Row 20
${ws}.Cells.Item(20, 4).NumberFormat = "@"
${ws}.Cells.Item(20, 4).Value = "1.003"
${ws}.Cells.Item(20, 5).NumberFormat = "@"
${ws}.Cells.Item(20, 5).Value = "  +0.06%  "

# Row 21
${ws}.Cells.Item(21, 5).NumberFormat = "@"
${ws}.Cells.Item(21, 5).Value = "  +1.90%  "

# Row 22
${ws}.Cells.Item(22, 5).NumberFormat = "@"
${ws}.Cells.Item(22, 5).Value = "  +1.63%  "

# Row 23
${ws}.Cells.Item(23, 4).NumberFormat = "@"
${ws}.Cells.Item(23, 4).Value = "6.096"
${ws}.Cells.Item(23, 5).NumberFormat = "@"
${ws}.Cells.Item(23, 5).Value = "  +2.37%  "

# Row 24
${ws}.Cells.Item(24, 5).NumberFormat = "@"
${ws}.Cells.Item(24, 5).Value = "  +0.05%  "

# Row 25
${ws}.Cells.Item(25, 4).NumberFormat = "@"
${ws}.Cells.Item(25, 4).Value = "148.23"
${ws}.Cells.Item(25, 5).NumberFormat = "@"
${ws}.Cells.Item(25, 5).Value = "  +2.18%  "

# Row 26
${ws}.Cells.Item(26, 5).NumberFormat = "@"
${ws}.Cells.Item(26, 5).Value = "  +12.15%  "

# Row 27
${ws}.Cells.Item(27, 4).NumberFormat = "@"
${ws}.Cells.Item(27, 4).Value = "0.1232"
${ws}.Cells.Item(27, 5).NumberFormat = "@"
${ws}.Cells.Item(27, 5).Value = "  +2.97%  "

# Row 28
${ws}.Cells.Item(28, 4).NumberFormat = "@"
${ws}.Cells.Item(28, 4).Value = "7.400"
${ws}.Cells.Item(28, 5).NumberFormat = "@"
${ws}.Cells.Item(28, 5).Value = "  +1.82%  "

# Row 29
${ws}.Cells.Item(29, 4).NumberFormat = "@"
${ws}.Cells.Item(29, 4).Value = "16.66"
${ws}.Cells.Item(29, 5).NumberFormat = "@"
${ws}.Cells.Item(29, 5).Value = "  +4.41%  "

# Row 30
${ws}.Cells.Item(30, 4).NumberFormat = "@"
${ws}.Cells.Item(30, 4).Value = "0.05537"
${ws}.Cells.Item(30, 5).NumberFormat = "@"
${ws}.Cells.Item(30, 5).Value = "  +1.53%  "

# Row 31
${ws}.Cells.Item(31, 5).NumberFormat = "@"
${ws}.Cells.Item(31, 5).Value = "  +2.53%  "

# Row 32
${ws}.Cells.Item(32, 5).NumberFormat = "@"
${ws}.Cells.Item(32, 5).Value = "  +2.31%  "

# Row 33
${ws}.Cells.Item(33, 4).NumberFormat = "@"
${ws}.Cells.Item(33, 4).Value = "3.462"
${ws}.Cells.Item(33, 5).NumberFormat = "@"
${ws}.Cells.Item(33, 5).Value = "  +2.76%  "

# Row 34
${ws}.Cells.Item(34, 5).NumberFormat = "@"
${ws}.Cells.Item(34, 5).Value = "  +6.24%  "

# Row 35
${ws}.Cells.Item(35, 4).NumberFormat = "@"
${ws}.Cells.Item(35, 4).Value = "0.9611"
${ws}.Cells.Item(35, 5).NumberFormat = "@"
${ws}.Cells.Item(35, 5).Value = "  +0.70%  "

# Row 36
${ws}.Cells.Item(36, 2).NumberFormat = "@"
${ws}.Cells.Item(36, 2).Value = "HuobiToken"
${ws}.Cells.Item(36, 3).NumberFormat = "@"
${ws}.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
${ws}.Cells.Item(36, 4).NumberFormat = "@"
${ws}.Cells.Item(36, 4).Value = "2.447"
${ws}.Cells.Item(36, 5).NumberFormat = "@"
${ws}.Cells.Item(36, 5).Value = "  +1.89%  "

# Row 37
${ws}.Cells.Item(37, 2).NumberFormat = "@"
${ws}.Cells.Item(37, 2).Value = "MXToken"
${ws}.Cells.Item(37, 3).NumberFormat = "@"
${ws}.Cells.Item(37, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
${ws}.Cells.Item(37, 4).NumberFormat = "@"
${ws}.Cells.Item(37, 4).Value = "2.821"
${ws}.Cells.Item(37, 5).NumberFormat = "@"
${ws}.Cells.Item(37, 5).Value = "  +1.49%  "

# Row 38
${ws}.Cells.Item(38, 4).NumberFormat = "@"
${ws}.Cells.Item(38, 4).Value = "0.5955"
${ws}.Cells.Item(38, 5).NumberFormat = "@"
${ws}.Cells.Item(38, 5).Value = "  +4.84%  "

# Row 39
${ws}.Cells.Item(39, 4).NumberFormat = "@"
${ws}.Cells.Item(39, 4).Value = "0.01647"
${ws}.Cells.Item(39, 5).NumberFormat = "@"
${ws}.Cells.Item(39, 5).Value = "  +3.83%  "

# Row 40
${ws}.Cells.Item(40, 4).NumberFormat = "@"
${ws}.Cells.Item(40, 4).Value = "5.929"
${ws}.Cells.Item(40, 5).NumberFormat = "@"
${ws}.Cells.Item(40, 5).Value = "  +0.90%  "

# Row 41
${ws}.Cells.Item(41, 4).NumberFormat = "@"
${ws}.Cells.Item(41, 4).Value = "1.058.99"
${ws}.Cells.Item(41, 5).NumberFormat = "@"
${ws}.Cells.Item(41, 5).Value = "  +3.01%  "

# Row 42
${ws}.Cells.Item(42, 4).NumberFormat = "@"
${ws}.Cells.Item(42, 4).Value = "0.8542"
${ws}.Cells.Item(42, 5).NumberFormat = "@"
${ws}.Cells.Item(42, 5).Value = "  +2.66%  "

# Row 43
${ws}.Cells.Item(43, 4).NumberFormat = "@"
${ws}.Cells.Item(43, 4).Value = "1.003"
${ws}.Cells.Item(43, 5).NumberFormat = "@"
${ws}.Cells.Item(43, 5).Value = "  +0.06%  "

# Row 44
${ws}.Cells.Item(44, 4).NumberFormat = "@"
${ws}.Cells.Item(44, 4).Value = "101.57"
${ws}.Cells.Item(44, 5).NumberFormat = "@"
${ws}.Cells.Item(44, 5).Value = "  +0.37%  "

# Row 45
${ws}.Cells.Item(45, 4).NumberFormat = "@"
${ws}.Cells.Item(45, 4).Value = "1.868.11"
${ws}.Cells.Item(45, 5).NumberFormat = "@"
${ws}.Cells.Item(45, 5).Value = "  +4.08%  "

# Row 46
${ws}.Cells.Item(46, 5).NumberFormat = "@"
${ws}.Cells.Item(46, 5).Value = "  +10.64%  "

# Row 47
${ws}.Cells.Item(47, 4).NumberFormat = "@"
${ws}.Cells.Item(47, 4).Value = "59.06"
${ws}.Cells.Item(47, 5).NumberFormat = "@"
${ws}.Cells.Item(47, 5).Value = "  +2.31%  "

# Row 48
${ws}.Cells.Item(48, 4).NumberFormat = "@"
${ws}.Cells.Item(48, 4).Value = "8.194"
${ws}.Cells.Item(48, 5).NumberFormat = "@"
${ws}.Cells.Item(48, 5).Value = "  +2.36%  "

# Row 49
${ws}.Cells.Item(49, 5).NumberFormat = "@"
${ws}.Cells.Item(49, 5).Value = "  +2.24%  "

# Row 50
${ws}.Cells.Item(50, 4).NumberFormat = "@"
${ws}.Cells.Item(50, 4).Value = "1.003"
${ws}.Cells.Item(50, 5).NumberFormat = "@"
${ws}.Cells.Item(50, 5).Value = "  +0.18%  "

# Row 51
${ws}.Cells.Item(51, 4).NumberFormat = "@"
${ws}.Cells.Item(51, 4).Value = "0.05277"
${ws}.Cells.Item(51, 5).NumberFormat = "@"
${ws}.Cells.Item(51, 5).Value = "  +1.54%  "
